$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the Timestamp column (A2:A97) forward by 31 days (one month later)
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value2 = $cell.Value2 + 31
}

# Update the Actual Production (MW) values for rows 21-47 to reflect the
# new day's solar production curve (Elnet added to Forecast Portfolio)
$newB = @(0, 2, 18, 45, 84, 139, 211, 300, 396, 508, 616, 748, 850, 950, 1033, 1101, 1166, 1241, 1304, 0, 0, 0, 0, 0, 0, 0, 0)

$r = 21
foreach ($val in $newB) {
    $ws.Cells.Item($r, 2).Value2 = $val
    $r++
}
